$wb = $excel.ActiveWorkbook

# ---- ALERTS sheet: append 3 row(s) starting at row 2 ----
$ws = $wb.Worksheets.Item("ALERTS")
$data = @()
$data += ,@("2026-01-28", "16:38:03", "16:00", "Bathroom", "MINIMAL", "MINIMAL ALERT: Bathroom occupied, no motion > 20s.")
$data += ,@("2026-01-28", "16:38:23", "16:00", "Bathroom", "MODERATE", "MODERATE ALERT: Bathroom occupied, no motion > 40s.")
$data += ,@("2026-01-28", "16:38:43", "16:00", "Bathroom", "CRITICAL", "CRITICAL ALERT: Bathroom occupied, no motion > 60s.")
$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).NumberFormat = "@"
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $row[$c]
    }
    $r++
}

# ---- PIR sheet: append 13 row(s) starting at row 34 ----
$ws = $wb.Worksheets.Item("PIR")
$data = @()
$data += ,@("2026-01-28", "16:37:52", "16:00", "Bathroom", "No Motion", "Inactive")
$data += ,@("2026-01-28", "16:37:52", "16:00", "Bathroom", "No Motion", "Inactive")
$data += ,@("2026-01-28", "16:37:57", "16:00", "Bathroom", "No Motion", "Inactive")
$data += ,@("2026-01-28", "16:38:02", "16:00", "Bathroom", "No Motion", "Inactive")
$data += ,@("2026-01-28", "16:38:07", "16:00", "Bathroom", "No Motion", "Inactive")
$data += ,@("2026-01-28", "16:38:12", "16:00", "Bathroom", "No Motion", "Inactive")
$data += ,@("2026-01-28", "16:38:17", "16:00", "Bathroom", "No Motion", "Inactive")
$data += ,@("2026-01-28", "16:38:22", "16:00", "Bathroom", "No Motion", "Inactive")
$data += ,@("2026-01-28", "16:38:27", "16:00", "Bathroom", "No Motion", "Inactive")
$data += ,@("2026-01-28", "16:38:32", "16:00", "Bathroom", "No Motion", "Inactive")
$data += ,@("2026-01-28", "16:38:37", "16:00", "Bathroom", "No Motion", "Inactive")
$data += ,@("2026-01-28", "16:38:42", "16:00", "Bathroom", "No Motion", "Inactive")
$data += ,@("2026-01-28", "16:38:47", "16:00", "Bathroom", "No Motion", "Inactive")
$r = 34
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).NumberFormat = "@"
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $row[$c]
    }
    $r++
}

# ---- Humidity sheet: append 14 row(s) starting at row 34 ----
$ws = $wb.Worksheets.Item("Humidity")
$data = @()
$data += ,@("2026-01-28", "16:37:52", "16:00", "Bathroom", "87.8%", "Active")
$data += ,@("2026-01-28", "16:37:54", "16:00", "Bathroom", "86.9%", "Active")
$data += ,@("2026-01-28", "16:37:58", "16:00", "Bathroom", "87.8%", "Active")
$data += ,@("2026-01-28", "16:38:02", "16:00", "Bathroom", "86.9%", "Active")
$data += ,@("2026-01-28", "16:38:06", "16:00", "Bathroom", "87.8%", "Active")
$data += ,@("2026-01-28", "16:38:10", "16:00", "Bathroom", "87.9%", "Active")
$data += ,@("2026-01-28", "16:38:14", "16:00", "Bathroom", "86.9%", "Active")
$data += ,@("2026-01-28", "16:38:22", "16:00", "Bathroom", "86.9%", "Active")
$data += ,@("2026-01-28", "16:38:26", "16:00", "Bathroom", "87.9%", "Active")
$data += ,@("2026-01-28", "16:38:30", "16:00", "Bathroom", "87.9%", "Active")
$data += ,@("2026-01-28", "16:38:35", "16:00", "Bathroom", "87.0%", "Active")
$data += ,@("2026-01-28", "16:38:39", "16:00", "Bathroom", "87.9%", "Active")
$data += ,@("2026-01-28", "16:38:47", "16:00", "Bathroom", "86.9%", "Active")
$data += ,@("2026-01-28", "16:38:51", "16:00", "Bathroom", "87.9%", "Active")
$r = 34
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).NumberFormat = "@"
    $ws.Cells.Item($r, 5).NumberFormat = "@"
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $row[$c]
    }
    $r++
}

# ---- Temperature sheet: append 14 row(s) starting at row 34 ----
$ws = $wb.Worksheets.Item("Temperature")
$data = @()
$data += ,@("2026-01-28", "16:37:52", "16:00", "Bathroom", "22.8C", "Active")
$data += ,@("2026-01-28", "16:37:55", "16:00", "Bathroom", "22.8C", "Active")
$data += ,@("2026-01-28", "16:37:58", "16:00", "Bathroom", "22.8C", "Active")
$data += ,@("2026-01-28", "16:38:02", "16:00", "Bathroom", "22.8C", "Active")
$data += ,@("2026-01-28", "16:38:06", "16:00", "Bathroom", "22.8C", "Active")
$data += ,@("2026-01-28", "16:38:11", "16:00", "Bathroom", "22.9C", "Active")
$data += ,@("2026-01-28", "16:38:15", "16:00", "Bathroom", "22.8C", "Active")
$data += ,@("2026-01-28", "16:38:23", "16:00", "Bathroom", "22.8C", "Active")
$data += ,@("2026-01-28", "16:38:27", "16:00", "Bathroom", "22.8C", "Active")
$data += ,@("2026-01-28", "16:38:31", "16:00", "Bathroom", "22.8C", "Active")
$data += ,@("2026-01-28", "16:38:35", "16:00", "Bathroom", "22.8C", "Active")
$data += ,@("2026-01-28", "16:38:39", "16:00", "Bathroom", "22.8C", "Active")
$data += ,@("2026-01-28", "16:38:47", "16:00", "Bathroom", "22.8C", "Active")
$data += ,@("2026-01-28", "16:38:51", "16:00", "Bathroom", "22.8C", "Active")
$r = 34
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).NumberFormat = "@"
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $row[$c]
    }
    $r++
}

Write-Output "Edit complete"